$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-9 (columns D, M, N, O, P, Q, S), resulting from
# the diff which reorders the price observations by date.

$data = @(
    @{Row=2; D=44176; M=250; N=7000;  O=7000;  P=7000;  Q="`$/caja 14 kilos empedrada"; S=500},
    @{Row=3; D=44397; M=60;  N=11000; O=11000; P=11000; Q="`$/caja 14 kilos";           S=786},
    @{Row=4; D=44309; M=300; N=7000;  O=7000;  P=7000;  Q="`$/caja 14 kilos empedrada"; S=500},
    @{Row=5; D=44400; M=100; N=10000; O=10000; P=10000; Q="`$/caja 14 kilos";           S=714},
    @{Row=6; D=44208; M=210; N=10000; O=10000; P=10000; Q="`$/caja 14 kilos empedrada"; S=714},
    @{Row=7; D=44351; M=300; N=10000; O=10000; P=10000; Q="`$/caja 14 kilos empedrada"; S=714},
    @{Row=8; D=44491; M=180; N=9000;  O=9000;  P=9000;  Q="`$/caja 14 kilos empedrada"; S=643},
    @{Row=9; D=44162; M=120; N=7000;  O=7000;  P=7000;  Q="`$/caja 14 kilos empedrada"; S=500}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # Column D - Fecha
    $ws.Cells.Item($r, 13).Value = $item.M   # Column M - Volumen
    $ws.Cells.Item($r, 14).Value = $item.N   # Column N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O   # Column O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P   # Column P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $item.Q   # Column Q - Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $item.S   # Column S - Precio $/Kg
}
